$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 21 (shifts existing rows 21-84 down to 24-87)
$ws.Rows("21:23").Insert()

# Row 21: new weekly price record
$ws.Range("A21").Value = 1
$ws.Range("B21").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C21").Value = 'Arica y Parinacota'
$ws.Range("D21").Value = '2023-01-26'
$ws.Range("E21").Value = 15
$ws.Range("F21").Value = 'Fruta'
$ws.Range("G21").Value = 100103
$ws.Range("H21").Value = 'Frutos de hueso (carozo)'
$ws.Range("I21").Value = 100103006
$ws.Range("J21").Value = 'Nectarín'
$ws.Range("K21").Value = 'Candy White'
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 300
$ws.Range("N21").Value = 24000
$ws.Range("O21").Value = 25000
$ws.Range("P21").Value = 24500
$ws.Range("Q21").Value = '$/bandeja 18 kilos granel'
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 1361
$ws.Range("T21").Value = 18

# Row 22: new weekly price record
$ws.Range("A22").Value = 1
$ws.Range("B22").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C22").Value = 'Arica y Parinacota'
$ws.Range("D22").Value = '2023-01-26'
$ws.Range("E22").Value = 15
$ws.Range("F22").Value = 'Fruta'
$ws.Range("G22").Value = 100103
$ws.Range("H22").Value = 'Frutos de hueso (carozo)'
$ws.Range("I22").Value = 100103006
$ws.Range("J22").Value = 'Nectarín'
$ws.Range("K22").Value = 'June Pearl'
$ws.Range("L22").Value = 'Segunda'
$ws.Range("M22").Value = 250
$ws.Range("N22").Value = 24000
$ws.Range("O22").Value = 25000
$ws.Range("P22").Value = 24500
$ws.Range("Q22").Value = '$/bandeja 18 kilos granel'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 1361
$ws.Range("T22").Value = 18

# Row 23: new weekly price record
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C23").Value = 'Arica y Parinacota'
$ws.Range("D23").Value = '2023-01-26'
$ws.Range("E23").Value = 15
$ws.Range("F23").Value = 'Fruta'
$ws.Range("G23").Value = 100103
$ws.Range("H23").Value = 'Frutos de hueso (carozo)'
$ws.Range("I23").Value = 100103006
$ws.Range("J23").Value = 'Nectarín'
$ws.Range("K23").Value = 'Venus'
$ws.Range("L23").Value = 'Segunda'
$ws.Range("M23").Value = 250
$ws.Range("N23").Value = 24000
$ws.Range("O23").Value = 25000
$ws.Range("P23").Value = 24500
$ws.Range("Q23").Value = '$/bandeja 18 kilos granel'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 1361
$ws.Range("T23").Value = 18

